{"js": "// Update the heading date/day-of-week line and every arithmetic problem in\n// the 20x5 table so the document content matches what was produced at\n// commit c986bee (an \"ADC/SBB within 100\" worksheet refresh).\n\n// index 0 = the heading paragraph above the table (the date line);\n// indices 1..100 = the 100 table-cell paragraphs, in row-major order\n// (row 1 col 1..5, row 2 col 1..5, ...).\nconst dateHeading = \"2024-11-16 Saturday\";\nconst cellValues = [\n  \"16+27=\", \"82-49=\", \"63-16=\", \"41-36=\", \"38+23=\", \"39+6=\", \"86-37=\", \"25+28=\",\n  \"90-25=\", \"13+19=\", \"91-19=\", \"93-16=\", \"98-79=\", \"32-27=\", \"91-3=\", \"19+53=\",\n  \"59+29=\", \"85+6=\", \"34-25=\", \"91-35=\", \"44+18=\", \"37+6=\", \"8+47=\", \"4+49=\",\n  \"42-33=\", \"65+19=\", \"18+25=\", \"36+17=\", \"51-44=\", \"38+36=\", \"48+49=\", \"5+49=\",\n  \"40-27=\", \"9+48=\", \"19+53=\", \"47+28=\", \"83-14=\", \"6+27=\", \"28+25=\", \"83-29=\",\n  \"25-9=\", \"57-8=\", \"22-14=\", \"52-28=\", \"72-34=\", \"22-9=\", \"7+89=\", \"65-48=\",\n  \"76-57=\", \"74-36=\", \"80-53=\", \"61-28=\", \"90-14=\", \"62-37=\", \"57+18=\", \"66+19=\",\n  \"64-25=\", \"90-41=\", \"71-32=\", \"60-57=\", \"17+55=\", \"19+25=\", \"53-36=\", \"8+83=\",\n  \"49+19=\", \"90-25=\", \"56+39=\", \"34-18=\", \"37+56=\", \"58+34=\", \"92-26=\", \"29+58=\",\n  \"84-46=\", \"6+7=\", \"31-26=\", \"66+25=\", \"16+9=\", \"72-4=\", \"18+46=\", \"80-23=\",\n  \"80-25=\", \"8+55=\", \"27+15=\", \"93-9=\", \"54+28=\", \"8+33=\", \"92-44=\", \"98-89=\",\n  \"69+24=\", \"40-16=\", \"93-55=\", \"43+29=\", \"20-6=\", \"38-9=\", \"17+74=\", \"85-79=\",\n  \"86-79=\", \"84-67=\", \"44-38=\", \"6+49=\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst expectedCount = 1 + cellValues.length;\nif (paragraphs.items.length !== expectedCount) {\n  throw new Error(\n    `Expected ${expectedCount} paragraphs (1 heading + ${cellValues.length} table cells), ` +\n    `found ${paragraphs.items.length}.`\n  );\n}\n\n// Paragraph 0: the date/day-of-week heading.\nparagraphs.items[0].getRange().insertText(dateHeading, Word.InsertLocation.replace);\n\n// Paragraphs 1..100: the table cells, in the same (row-major) order they\n// appear in the document.\nfor (let i = 0; i < cellValues.length; i++) {\n  paragraphs.items[i + 1].getRange().insertText(cellValues[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the heading date/day-of-week line and every arithmetic problem in\n# the 20x5 table so the document content matches what was produced at\n# commit c986bee (an \"ADC/SBB within 100\" worksheet refresh).\n\n$d = $word.ActiveDocument\n\n# --- 1. Heading paragraph (the date line shown above the table) -----------\n$dateHeading = '2024-11-16 Saturday'\n$d.Paragraphs.Item(1).Range.Text = $dateHeading\n\n# --- 2. Table cells (20 rows x 5 columns = 100 problems), listed in\n#        row-major order (row 1 col 1..5, row 2 col 1..5, ...) to match the\n#        order the cells appear in the document.\n$cellValues = @(\n    '16+27=', '82-49=', '63-16=', '41-36=', '38+23=', '39+6=', '86-37=', '25+28=',\n    '90-25=', '13+19=', '91-19=', '93-16=', '98-79=', '32-27=', '91-3=', '19+53=',\n    '59+29=', '85+6=', '34-25=', '91-35=', '44+18=', '37+6=', '8+47=', '4+49=',\n    '42-33=', '65+19=', '18+25=', '36+17=', '51-44=', '38+36=', '48+49=', '5+49=',\n    '40-27=', '9+48=', '19+53=', '47+28=', '83-14=', '6+27=', '28+25=', '83-29=',\n    '25-9=', '57-8=', '22-14=', '52-28=', '72-34=', '22-9=', '7+89=', '65-48=',\n    '76-57=', '74-36=', '80-53=', '61-28=', '90-14=', '62-37=', '57+18=', '66+19=',\n    '64-25=', '90-41=', '71-32=', '60-57=', '17+55=', '19+25=', '53-36=', '8+83=',\n    '49+19=', '90-25=', '56+39=', '34-18=', '37+56=', '58+34=', '92-26=', '29+58=',\n    '84-46=', '6+7=', '31-26=', '66+25=', '16+9=', '72-4=', '18+46=', '80-23=',\n    '80-25=', '8+55=', '27+15=', '93-9=', '54+28=', '8+33=', '92-44=', '98-89=',\n    '69+24=', '40-16=', '93-55=', '43+29=', '20-6=', '38-9=', '17+74=', '85-79=',\n    '86-79=', '84-67=', '44-38=', '6+49='\n)\n\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\nif (($rows * $cols) -ne $cellValues.Count) {\n    throw \"Expected $($cellValues.Count) cells, table has $rows x $cols = $($rows * $cols).\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $t.Cell($r, $c).Range.Text = $cellValues[$i]\n        $i++\n    }\n}\n"}
